$d = $word.ActiveDocument

# --- 1) Resize the first table's grid columns (in points; twips/20) ---
$tbl = $d.Tables(1)
$tbl.Columns(1).Width = 800 / 20
$tbl.Columns(2).Width = 560 / 20
$tbl.Columns(3).Width = 6560 / 20

# --- 2) Update the text in problem 4's solution cell ---
$d.Content.Find.Execute(
    "A discrete random variable is something that varies following a specific pattern",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A probability of 1 implies an event is certain to happen. A probability of 0",
    2)

$d.Content.Find.Execute(
    "or distribution over the long run. They are discrete if they can be listed.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "implies it is impossible to happen, or certain to not happen.",
    2)
